# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Each entry maps a row number to the new value that should be written
# into column F of that row.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览" = @{
        2  = 626
        4  = 1301
        5  = 1171
        6  = 14359
        7  = 16630
        9  = 104
        10 = 18
        19 = 106
        21 = 1267
        24 = 41
        25 = 20
        27 = 6765
        29 = 22
        30 = 1123
        33 = 5763
        34 = 108
        36 = 193
        37 = 4844
        38 = 19
    }
    "全部类型" = @{
        2  = 626
        4  = 1301
        5  = 1171
        6  = 14359
        7  = 16630
        9  = 104
        10 = 18
        19 = 106
        21 = 1267
        25 = 41
        26 = 20
        28 = 6765
        30 = 22
        31 = 1123
        36 = 5763
        37 = 108
        39 = 193
        40 = 4844
        41 = 19
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $sheetUpdates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
